$d = $word.ActiveDocument
$xml = $d.WordOpenXML

function Remove-DuplicateStyleBlocks {
    param(
        [string]$Text,
        [string]$Pattern
    )
    $re = [regex]$Pattern
    $matches = $re.Matches($Text)
    if ($matches.Count -le 1) {
        return $Text
    }
    $result = ""
    $lastEnd = 0
    for ($i = 0; $i -lt $matches.Count; $i++) {
        $m = $matches[$i]
        $result += $Text.Substring($lastEnd, $m.Index - $lastEnd)
        if ($i -eq 0) {
            $result += $m.Value
        }
        $lastEnd = $m.Index + $m.Length
    }
    $result += $Text.Substring($lastEnd)
    return $result
}

$xml = Remove-DuplicateStyleBlocks -Text $xml -Pattern '<w:style w:type="table" w:styleId="Grilledutableau">.*?</w:style>'
$xml = Remove-DuplicateStyleBlocks -Text $xml -Pattern '<w:style w:type="table" w:default="1" w:styleId="TableauNormal">.*?</w:style>'

$d.WordOpenXML = $xml
